$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.335.54'
$ws.Range('E2').Value = '  -0.71%  '
$ws.Range('D3').Value = '3.216.94'
$ws.Range('E3').Value = '  +0.53%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '606.45'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.61'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.51%  '
$ws.Range('D8').Value = '3.215.70'
$ws.Range('E8').Value = '  +0.54%  '
$ws.Range('E9').Value = '  -2.06%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.70'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -3.94%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.500'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -3.16%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000268'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '38.28'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.41%  '
$ws.Range('D15').Value = '3.744.59'
$ws.Range('E15').Value = '  +0.54%  '
$ws.Range('D16').Value = '66.413.73'
$ws.Range('E16').Value = '  -0.50%  '
$ws.Range('D17').Value = '3.217.36'
$ws.Range('E17').Value = '  +0.52%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.24'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -3.14%  '
$ws.Range('E19').Value = '  +0.93%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '506.14'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -3.03%  '
$ws.Range('E21').Value = '  -1.69%  '
$ws.Range('E22').Value = '  -2.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.95'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -3.18%  '
$ws.Range('E24').Value = '  -3.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.13'
$ws.Range('D25').ClearFormats()
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.156'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +72.82%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('E28').Value = '  -0.95%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.00'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -2.94%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.36'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.34%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.88'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -4.19%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.89'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.95%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '28.17'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.67%  '
$ws.Range('E34').Value = '  +0.12%  '
$ws.Range('E35').Value = '  -5.23%  '
$ws.Range('E36').Value = '  -2.80%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '55.31'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '499.81'
$ws.Range('D38').ClearFormats()
$ws.Range('D39').Value = '0.0₃0768'
$ws.Range('E39').Value = '  +11.34%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0416'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -2.48%  '
$ws.Range('E41').Value = '  +1.69%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.01'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +3.66%  '
$ws.Range('E43').Value = '  -2.38%  '
$ws.Range('E44').Value = '  -2.97%  '
$ws.Range('D45').Value = '2.923.71'
$ws.Range('E45').Value = '  +0.92%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.45'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.98%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '28.04'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.16%  '
$ws.Range('E48').Value = '  +1.65%  '
$ws.Range('E49').Value = '  -0.69%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '121.23'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.09%  '
